$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original sheet and add the two new sheets (at the end, in
#    order), matching the target sheet order: CPC, CS, CATCH.
# ---------------------------------------------------------------------------
$cpc = $wb.Worksheets.Item(1)
$cpc.Name = "CPC"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cs = $wb.Worksheets.Add($null, $lastSheet)
$cs.Name = "CS"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$catch = $wb.Worksheets.Add($null, $lastSheet2)
$catch.Name = "CATCH"

# ---------------------------------------------------------------------------
# 2. CPC sheet: fill in column C (Example / Format) for every field, with a
#    centered alignment style, and a percentage for the last (numeric) row.
# ---------------------------------------------------------------------------
$cpc.Range("C2").Value = "MDV"
$cpc.Range("C2").HorizontalAlignment = -4108

$cpc.Range("C3").Value = "Maldives"
$cpc.Range("C3").HorizontalAlignment = -4108

$cpc.Range("C4").Value = "Maldives"
$cpc.Range("C4").HorizontalAlignment = -4108

$cpc.Range("C5").Value = "CP"
$cpc.Range("C5").HorizontalAlignment = -4108

$cpc.Range("C6").Value = "Contracting Party"
$cpc.Range("C6").HorizontalAlignment = -4108

$cpc.Range("C7").Value = "true or false"
$cpc.Range("C7").HorizontalAlignment = -4108

$cpc.Range("C8").Value = "true or false"
$cpc.Range("C8").HorizontalAlignment = -4108

$cpc.Range("C9").Value = "true or false"
$cpc.Range("C9").HorizontalAlignment = -4108

$cpc.Range("C10").Value = "916,244 km2"
$cpc.Range("C10").HorizontalAlignment = -4108

$cpc.Range("C11").Value = 0.0149
$cpc.Range("C11").NumberFormat = "0.00%"
$cpc.Range("C11").HorizontalAlignment = -4108

# Selection moves from C7 to B7, and CPC is no longer the selected tab.
$cpc.Range("B7").Select()

# ---------------------------------------------------------------------------
# 3. CS sheet: headers, field metadata, and "Used in Option" column.
# ---------------------------------------------------------------------------
$cs.Range("A1").Value = "Field Name"
$cs.Range("B1").Value = "Description"
$cs.Range("C1").Value = "Used in Option"

$cs.Range("A2").Value = "CODE"
$cs.Range("B2").Value = "Mnemonic code, generally the ISO3 code"
$cs.Range("C2").Value = "-"

$cs.Range("A3").Value = "NAME_EN"
$cs.Range("B3").Value = "Official English Name"
$cs.Range("C3").Value = "-"

$cs.Range("A4").Value = "NAME_FR"
$cs.Range("B4").Value = "Official French Name"
$cs.Range("C4").Value = "-"

$cs.Range("A5").Value = "DEVELOPMENT_STATUS"
$cs.Range("B5").Value = "Development classification: Least Developed (LD), Developing (DG), Developed (DE)"
$cs.Range("C5").Value = "-"

$cs.Range("A6").Value = "PER_CAPITA_FISH_CONSUMPTION_KG"
$cs.Range("B6").Value = "Per capita fish consumption (kg/person/year)"
$cs.Range("C6").Value = "Option 1"

$cs.Range("A7").Value = "CUV_INDEX"
$cs.Range("B7").Value = "Commonwealth Universal Vulnerability Index"
$cs.Range("C7").Value = "Option 1"

$cs.Range("A8").Value = "PROP_WORKERS_EMPLOYED_SSF"
$cs.Range("B8").Value = "Percentage of fish workers employed in small-scale and artisanal fisheries"
$cs.Range("C8").Value = "Option 1"

$cs.Range("A9").Value = "SIDS_STATUS"
$cs.Range("B9").Value = "Whether the CPC is a Small Island Developing State"
$cs.Range("C9").Value = "Both"

$cs.Range("A10").Value = "PROP_FISHERIES_CONTRIBUTION_GDP"
$cs.Range("B10").Value = "Percentage contribution of fisheries to Gross Domestic Product"
$cs.Range("C10").Value = "Option 1"

$cs.Range("A11").Value = "PROP_EXPORT_VALUE_FISHERY"
$cs.Range("B11").Value = "Percentage of total export value contributed by fisheries"
$cs.Range("C11").Value = "Option 1"

$cs.Range("A12").Value = "HDI_STATUS"
$cs.Range("B12").Value = "Human Development Index status"
$cs.Range("C12").Value = "Option 2"

$cs.Range("A13").Value = "GNI_STATUS"
$cs.Range("B13").Value = "Gross National Income status"
$cs.Range("C13").Value = "Option 2"

# Column C is formatted as Text for the whole used block (rows 2-14).
$cs.Range("C2:C14").NumberFormat = "@"

# B9 and B11 reuse the existing wrap/vertical-center style (same as used on
# the CPC sheet) - copy format only so no new cellXf is minted.
$cpc.Range("B7").Copy()
$cs.Range("B9").PasteSpecial(-4122)
$cs.Range("B11").PasteSpecial(-4122)

# Two trailing formatted-but-empty cells.
$cs.Range("B16").PasteSpecial(-4122)
$cs.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths (best-effort match of the recorded bestFit widths).
$cs.Columns.Item(1).ColumnWidth = 33.666666666666664
$cs.Columns.Item(2).ColumnWidth = 74.83333333333333
$cs.Columns.Item(3).ColumnWidth = 22.666666666666668

$cs.Range("C11").Select()
$cs.Activate()
